# Update a batch of numeric values in Sheet1 as per the commit's data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value  = 13.111
$ws.Range("A4").Value  = -21.582
$ws.Range("D4").Value  = -8.134
$ws.Range("E4").Value  = 13.019
$ws.Range("D5").Value  = -8.581999999999999
$ws.Range("A6").Value  = -21.108
$ws.Range("A7").Value  = -21.047
$ws.Range("D8").Value  = -8.204000000000001
$ws.Range("E9").Value  = 12.961
$ws.Range("E11").Value = 12.852
$ws.Range("E14").Value = 13.06
$ws.Range("A16").Value = -20.667
$ws.Range("D16").Value = -8.600999999999999
$ws.Range("E18").Value = 12.596
$ws.Range("A20").Value = -22.035
$ws.Range("D22").Value = -8.134
$ws.Range("E25").Value = 12.791
